# Add a new row 12 ("2021年") to Sheet1, continuing the yearly series that
# currently ends at row 11 ("2020年"), covering columns A:BK.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 12
$templateRow = 11

# Copy the formatting of the template row's label cell (A11) onto the new
# label cell (A12) so it keeps the same style (bold, centered, bordered).
$ws.Range("A$templateRow").Copy() | Out-Null
$ws.Range("A$newRow").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Year label for the new row.
$ws.Range("A$newRow").Value = "2021年"

# Numeric data for the new row, column letter -> value.
$values = @{
    "B"  = 198;
    "C"  = 1848;
    "D"  = 1354;
    "E"  = 20;
    "F"  = 68;
    "G"  = 114;
    "H"  = 4674;
    "I"  = 1077;
    "J"  = 5425;
    "L"  = 227;
    "M"  = 1420;
    "N"  = 229;
    "P"  = 641;
    "Q"  = 254;
    "R"  = 263;
    "S"  = 727;
    "T"  = 175;
    "U"  = 390;
    "V"  = 2058;
    "W"  = 32;
    "X"  = 2964;
    "Y"  = 4412;
    "Z"  = 5275;
    "AA" = 1386;
    "AC" = 194;
    "AD" = 94;
    "AE" = 235;
    "AF" = 4445;
    "AG" = 5636;
    "AH" = 971;
    "AI" = 59;
    "AJ" = 280;
    "AK" = 325;
    "AN" = 504;
    "AO" = 484;
    "AP" = 1049;
    "AR" = 72;
    "AS" = 1818;
    "AT" = 139;
    "AU" = 5;
    "AV" = 399;
    "AW" = 3599;
    "AX" = 1088;
    "AY" = 1465;
    "AZ" = 237;
    "BB" = 595;
    "BC" = 3338;
    "BD" = 664;
    "BE" = 7951;
    "BG" = 72526;
    "BH" = 378;
    "BJ" = 890;
    "BK" = 381;
}

foreach ($col in $values.Keys) {
    $ws.Range("$col$newRow").Value = $values[$col]
}

# These columns (same countries that have no value in the template row
# either) stay blank in the new row - nothing to write for them.

Write-Host "Row $newRow populated."
